$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) previously carried a custom "applyAlignment" style
# (wrap-text look) on every cell. Drop that formatting so the header cells
# fall back to the workbook's default style.
$ws.Rows(1).ClearFormats()

# Column E (Etagenhoehe / height_of_floors) is a numeric column formatted
# with a 2-decimal number format ("0.00"); restore that on the header cell
# E1 (it lost it when we cleared the row's formatting above) and make sure
# every E-column value, including the new row below, uses it too.
$ws.Range("E1").NumberFormat = "0.00"

# Append the new building record as row 6.
$ws.Range("B6").Value = "Griessgram"
$ws.Range("C6").Value = 1998
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 4.76
$ws.Range("E6").NumberFormat = "0.00"
$ws.Range("F6").Value = 5000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0

# Mirror the author's final selection.
$null = $ws.Range("K6").Select()
